$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2025-05")

# Update existing row 6 ("hp mini 400 G9" / "i7-13700T") with new benchmark results
$ws.Cells.Item(6, 3).Value = 4826
$ws.Cells.Item(6, 4).Value = 235712000000

# Add new row 9 with the older result for the same CPU, now annotated with a source link
$ws.Cells.Item(9, 2).Value = "i7-13700T"
$ws.Cells.Item(9, 3).Value = 1840
$ws.Cells.Item(9, 4).Value = 29037000000
$ws.Cells.Item(9, 4).NumberFormat = $ws.Cells.Item(6, 4).NumberFormat

$link = $ws.Hyperlinks.Add($ws.Cells.Item(9, 6), "https://openbenchmarking.org/result/2504309-SAIH-KREIER840 ")

# Update the active selection on the sheet
$ws.Range("E12").Select() | Out-Null
